# Build site at 2022-09-26 16:07:08 UTC
# Re-lays out the "Docentes responsaveis:" / "Programa resumido:" / "Programa:" /
# "Avaliacao:" block of the LOB1238 syllabus sheet: the old row 22 (long
# bibliography row) is dropped, and rows 10, 13-21 are rewritten to their new
# values, including a handful of row-height tweaks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The final sheet only spans to row 21 - the old last row (bibliography) goes away.
$ws.Rows.Item(22).Delete()

# Row 10: B/C value changes (keeps label "Objetivos:" in A10).
$ws.Range("B10").Value = "5840942 - Marco Aurélio Kondracki de Alcântara"
$ws.Range("C10").Value = "5840942 - Marco Aurélio Kondracki de Alcântara"

# Row 13: now carries a label in A13 plus new B/C value.
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Rows.Item(13).RowHeight = 60

# Row 14: label changes, B/C cleared (no longer carries a value).
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").ClearContents()
$ws.Range("C14").ClearContents()

# Row 15: label changes, B/C now carries the activation date value.
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2012"
$ws.Range("C15").Value = "01/01/2012"
$ws.Rows.Item(15).RowHeight = 120

# Row 16: label changes, B/C cleared.
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").ClearContents()
$ws.Range("C16").ClearContents()

# Row 17: label changes only.
$ws.Range("A17").Value = "Avaliação:"
$ws.Rows.Item(17).RowHeight = 15
$ws.Rows.Item(17).AutoFit() | Out-Null

# Row 18: label changes, B/C now carries the docente value.
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "5840942 - Marco Aurélio Kondracki de Alcântara"
$ws.Range("C18").Value = "5840942 - Marco Aurélio Kondracki de Alcântara"
$ws.Rows.Item(18).RowHeight = 60

# Row 19: label changes only (B/C keep "Aula expositiva e exercícios dirigidos.").
$ws.Range("A19").Value = "Critério:"

# Row 20: label changes only (B/C keep "Média ponderada de exercícios e provas.").
$ws.Range("A20").Value = "Norma de recuperação:"

# Row 21: label changes only (B/C keep "Prova única com nota igual ou superior a 5,0.").
$ws.Range("A21").Value = "Bibliografia:"
$ws.Rows.Item(21).RowHeight = 120
